$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings used in row 2: User Name, Exam Password values
$ws.Range("A2").Value = "test114"
$ws.Range("C2").Value = "narendra62"
$ws.Range("D2").Value = "T3#%gA2b"

# Update Candidate ID value
$ws.Range("B2").Value = 23071147
